$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 824, shifting existing rows
# 824-896 down to 826-898 (mirrors the XML diff: dimension grows from
# A1:R896 to A1:R898, and every row from 824 onward shifts down by 2).
$ws.Rows("824:825").Insert()

# Fill the two newly-inserted rows with their data.
# Row 824 (new)
$ws.Range("A824").Value = 10
$ws.Range("B824").Value = "Vega Modelo de Temuco"
$ws.Range("C824").Value = "La Araucanía"
$ws.Range("D824").Value = 45166
$ws.Range("E824").Value = 9
$ws.Range("F824").Value = 100112045
$ws.Range("G824").Value = "Zapallo"
$ws.Range("H824").Value = "Camote"
$ws.Range("I824").Value = "1a (guarda)"
$ws.Range("J824").Value = 900
$ws.Range("K824").Value = 900
$ws.Range("L824").Value = 1000
$ws.Range("M824").Value = 944
$ws.Range("N824").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O824").Value = "Región del Maule"
$ws.Range("P824").Value = 944
$ws.Range("Q824").Value = 1
$ws.Range("R824").Value = "Hortaliza"

# Row 825 (new)
$ws.Range("A825").Value = 10
$ws.Range("B825").Value = "Vega Modelo de Temuco"
$ws.Range("C825").Value = "La Araucanía"
$ws.Range("D825").Value = 45166
$ws.Range("E825").Value = 9
$ws.Range("F825").Value = 100112045
$ws.Range("G825").Value = "Zapallo"
$ws.Range("H825").Value = "Camote"
$ws.Range("I825").Value = "1a nueva(o)"
$ws.Range("J825").Value = 400
$ws.Range("K825").Value = 900
$ws.Range("L825").Value = 900
$ws.Range("M825").Value = 900
$ws.Range("N825").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O825").Value = "Perú"
$ws.Range("P825").Value = 900
$ws.Range("Q825").Value = 1
$ws.Range("R825").Value = "Hortaliza"
